$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = '30.261.63'
$r.ClearFormats()
$ws.Range("E2").Value = '  -0.45%  '
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = '1.858.06'
$r.ClearFormats()
$ws.Range("E3").Value = '  -1.06%  '
$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = '1.001'
$r.ClearFormats()
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("E5").Value = '  -2.22%  '
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '1.001'
$r.ClearFormats()
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("E7").Value = '  -0.58%  '
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = '0.2755'
$r.ClearFormats()
$ws.Range("E8").Value = '  -2.51%  '
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = '0.06417'
$r.ClearFormats()
$ws.Range("E9").Value = '  -1.62%  '
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = '1.813.22'
$r.ClearFormats()
$ws.Range("E10").Value = '  -3.41%  '
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = '0.07428'
$r.ClearFormats()
$ws.Range("E11").Value = '  -0.52%  '
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = '16.08'
$r.ClearFormats()
$ws.Range("E12").Value = '  -3.54%  '
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = '4.992'
$r.ClearFormats()
$ws.Range("E13").Value = '  -2.09%  '
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = '85.18'
$r.ClearFormats()
$ws.Range("E14").Value = '  -3.51%  '
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = '0.6322'
$r.ClearFormats()
$ws.Range("E15").Value = '  -4.33%  '
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = '30.219.53'
$r.ClearFormats()
$ws.Range("E16").Value = '  -0.53%  '
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = '1.001'
$r.ClearFormats()
$ws.Range("E17").Value = '  +0.07%  '
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = '12.79'
$r.ClearFormats()
$ws.Range("E18").Value = '  -3.90%  '
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = '0.000007312'
$r.ClearFormats()
$ws.Range("E19").Value = '  -3.95%  '
$ws.Range("B20").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C20").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '2.098.60'
$r.ClearFormats()
$ws.Range("E20").Value = '  -0.92%  '
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = '223.70'
$r.ClearFormats()
$ws.Range("E21").Value = '  +2.18%  '
$ws.Range("B22").Value = 'BinanceUSD'
$ws.Range("C22").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = '1.002'
$r.ClearFormats()
$ws.Range("E22").Value = '  +0.17%  '
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '5.099'
$r.ClearFormats()
$ws.Range("E23").Value = '  -3.85%  '
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = '5.994'
$r.ClearFormats()
$ws.Range("E24").Value = '  -3.52%  '
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = '166.85'
$r.ClearFormats()
$ws.Range("E25").Value = '  -0.45%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = '9.207'
$r.ClearFormats()
$ws.Range("E26").Value = '  -1.80%  '
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = '17.77'
$r.ClearFormats()
$ws.Range("E27").Value = '  -3.72%  '
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = '1.864'
$r.ClearFormats()
$ws.Range("E28").Value = '  -5.75%  '
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = '0.1032'
$r.ClearFormats()
$ws.Range("E29").Value = '  +9.90%  '
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = '1.377'
$r.ClearFormats()
$ws.Range("E30").Value = '  -5.85%  '
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = '4.212'
$r.ClearFormats()
$ws.Range("E31").Value = '  -2.54%  '
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = '3.891'
$r.ClearFormats()
$ws.Range("E32").Value = '  -3.67%  '
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = '0.04883'
$r.ClearFormats()
$ws.Range("E33").Value = '  -3.31%  '
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = '1.149'
$r.ClearFormats()
$ws.Range("E34").Value = '  -4.54%  '
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = '0.7269'
$r.ClearFormats()
$ws.Range("E35").Value = '  -2.68%  '
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = '1.000'
$r.ClearFormats()
$ws.Range("E36").Value = '  +0.21%  '
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = '2.680'
$r.ClearFormats()
$ws.Range("E37").Value = '  -1.17%  '
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '0.01896'
$r.ClearFormats()
$ws.Range("E38").Value = '  +3.99%  '
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = '2.625'
$r.ClearFormats()
$ws.Range("E39").Value = '  +0.42%  '
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = '0.9018'
$r.ClearFormats()
$ws.Range("E40").Value = '  -0.43%  '
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = '1.974'
$r.ClearFormats()
$ws.Range("E41").Value = '  -4.79%  '
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = '105.41'
$r.ClearFormats()
$ws.Range("E42").Value = '  -1.42%  '
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '0.9942'
$r.ClearFormats()
$ws.Range("E43").Value = '  -1.00%  '
$ws.Range("E44").Value = '  -4.35%  '
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '5.541'
$r.ClearFormats()
$ws.Range("E45").Value = '  -6.14%  '
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = '7.060'
$r.ClearFormats()
$ws.Range("E46").Value = '  -4.89%  '
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = '61.03'
$r.ClearFormats()
$ws.Range("E47").Value = '  -5.22%  '
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = '0.1203'
$r.ClearFormats()
$ws.Range("E48").Value = '  -6.13%  '
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = '8.780'
$r.ClearFormats()
$ws.Range("E49").Value = '  -1.44%  '
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = '1.402'
$r.ClearFormats()
$ws.Range("E50").Value = '  -5.45%  '
$ws.Range("B51").Value = 'Elrond'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = '32.98'
$r.ClearFormats()
$ws.Range("E51").Value = '  -2.25%  '
